$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("GLOBAL RESULTS")
$ws.Range("C3").Value = 12.191162919778115
$ws.Range("C5").Value = 0.11524794485395798
$ws.Range("C7").Value = 49.83250922882285
$ws.Range("C9").Value = 4.760625171259168
$ws.Range("C13").Value = 11.82398098424618
$ws.Range("C15").Value = 0.26437248358866916
$ws.Range("C17").Value = 34.6650755106766
$ws.Range("C19").Value = 10.9206138257423
$ws.Range("C23").Value = 11.82398098424618
$ws.Range("C25").Value = 0.26437248358866916
$ws.Range("C27").Value = 34.6650755106766
$ws.Range("C29").Value = 10.9206138257423
$ws.Range("C33").Value = 11.82398098424618
$ws.Range("C35").Value = 0.26437248358866916
$ws.Range("C37").Value = 34.6650755106766
$ws.Range("C39").Value = 10.9206138257423
$ws.Range("C43").Value = 12.002801810951546
$ws.Range("C45").Value = 0.17950228112805888
$ws.Range("C47").Value = 42.05174891155756
$ws.Range("C49").Value = 7.4148227017805155
$ws.Range("C53").Value = 11.989289361632132
$ws.Range("C55").Value = 0.44754789626086433
$ws.Range("C57").Value = 41.49358098555356
$ws.Range("C59").Value = 18.487165068179394
$ws.Range("C62").Value = 20.340187707059705
$ws.Range("C63").Value = 41.49358098555356
$ws.Range("C64").Value = 54.24025116750932
$ws.Range("C69").Value = 45821.44336009229
$ws.Range("C70").Value = 750311.9952829664
$ws.Range("C71").Value = 704490.5519228742
$ws.Range("C76").Value = 48619.36661211877
$ws = $wb.Worksheets.Item("LANDING GEARS")
$ws.Range("C5").Value = -18.416904856045285
$ws.Range("C6").Value = -18.416904856045317
$ws.Range("C7").Value = 10.538941955553913
$ws.Range("C8").Value = 10.538941955553913
$ws.Range("C9").Value = 10.538941955553913
$ws.Range("C10").Value = 10.538941955553913
$ws.Range("C23").Value = 10.538941955553913
